$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 14 with the new "Cache I/D 16k" test entry
$ws.Range("B14").Value = "Cache I/D 16k"

# Setting .Value directly on C14 would drop its existing quote-prefixed
# center-aligned format, so set the value then restore the original
# format (copied from the identically-styled C13 cell) via PasteSpecial.
$ws.Range("C14").Value = "-O3"
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122)

$ws.Range("D14").Value = 32
$ws.Range("E14").Value = 108

# Update the selection to reflect the new active cell
$ws.Range("F16").Select()
